$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.968.54'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.574.68'
$ws.Range('E3').Value = '  +2.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.22%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.29'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0813'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '2.972.54'
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.643.52'
$ws.Range('E15').Value = '  +1.98%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.844'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.13%  '
$ws.Range('D18').Value = '43.020.47'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.72'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').Value = '0.0₃0970'
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.91%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '40.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.97%  '
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('E31').Value = '  -1.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '157.78'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.43'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.48%  '
$ws.Range('E34').Value = '  -1.22%  '
$ws.Range('E35').Value = '  +3.07%  '
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.79%  '
$ws.Range('E38').Value = '  +10.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.113'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.67'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +8.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0304'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.20%  '
$ws.Range('D46').Value = '2.002.83'
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.93'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').Value = '2.823.24'
$ws.Range('E48').Value = '  +2.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.198'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.13'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.76%  '
